$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 185.77777
$ws.Range("I33").Value = 177.4762
$ws.Range("J33").Value = 214.83333
$ws.Range("K33").Value = 177.4762
$ws.Range("L33").Value = 214.83333
$ws.Range("M33").Value = 51.52379999999999
$ws.Range("N33").Value = -672.8333299999999
$ws.Range("H62").Value = 9271.143
$ws.Range("I62").Value = 7483
$ws.Range("J62").Value = 20000
$ws.Range("K62").Value = 7483
$ws.Range("L62").Value = 20000
$ws.Range("M62").Value = -6859
$ws.Range("N62").Value = -21248
$ws.Range("H64").Value = 3218.75
$ws.Range("J64").Value = 3387.5
$ws.Range("L64").Value = 3387.5
$ws.Range("N64").Value = -3883.5
$ws.Range("H65").Value = 9271.143
$ws.Range("I65").Value = 7483
$ws.Range("J65").Value = 20000
$ws.Range("K65").Value = 37415
$ws.Range("L65").Value = 100000
$ws.Range("M65").Value = -34295
$ws.Range("N65").Value = -106240
$ws.Range("H67").Value = 3218.75
$ws.Range("J67").Value = 3387.5
$ws.Range("L67").Value = 3387.5
$ws.Range("N67").Value = -5103.5
$ws.Range("H76").Value = 10108488
$ws.Range("I76").Value = 8896.471
$ws.Range("J76").Value = 20839304
$ws.Range("K76").Value = 8896.471
$ws.Range("L76").Value = 20839304
$ws.Range("M76").Value = -8581.471
$ws.Range("N76").Value = -20839934
$ws.Range("H79").Value = 10108488
$ws.Range("I79").Value = 8896.471
$ws.Range("J79").Value = 20839304
$ws.Range("K79").Value = 8896.471
$ws.Range("L79").Value = 20839304
$ws.Range("M79").Value = -7804.471
$ws.Range("N79").Value = -20841488
$ws.Range("H112").Value = 400
$ws.Range("I112").Value = 400
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 1200
$ws.Range("L112").Value = 0
$ws.Range("M112").Value = -92
$ws.Range("N112").ClearContents()
$ws.Range("H125").Value = 3570
$ws.Range("I125").Value = 1910
$ws.Range("J125").Value = 4400
$ws.Range("K125").Value = 17190
$ws.Range("L125").Value = 39600
$ws.Range("M125").Value = -14730
$ws.Range("N125").Value = -44520

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2552.75
$ws.Range("I63").Value = 2377.6667
$ws.Range("J63").Value = 2611.111
$ws.Range("K63").Value = 2377.6667
$ws.Range("L63").Value = 2611.111
$ws.Range("M63").Value = -1691.6667
$ws.Range("N63").Value = -3983.111
$ws.Range("H66").Value = 2552.75
$ws.Range("I66").Value = 2377.6667
$ws.Range("J66").Value = 2611.111
$ws.Range("K66").Value = 11888.3335
$ws.Range("L66").Value = 13055.555
$ws.Range("M66").Value = -8456.333500000001
$ws.Range("N66").Value = -19919.555
$ws.Range("H102").Value = 1229.4231
$ws.Range("I102").Value = 1229.7273
$ws.Range("J102").Value = 1227.75
$ws.Range("K102").Value = 1229.7273
$ws.Range("L102").Value = 1227.75
$ws.Range("M102").Value = 392.2727
$ws.Range("N102").Value = -4471.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1226582.4
$ws.Range("I86").Value = 3025.6
$ws.Range("J86").Value = 2586090
$ws.Range("K86").Value = 3025.6
$ws.Range("L86").Value = 2586090
$ws.Range("M86").Value = -1902.6
$ws.Range("N86").Value = -2588336
$ws.Range("H89").Value = 1226582.4
$ws.Range("I89").Value = 3025.6
$ws.Range("J89").Value = 2586090
$ws.Range("K89").Value = 15128
$ws.Range("L89").Value = 12930450
$ws.Range("M89").Value = -9512
$ws.Range("N89").Value = -12941682
$ws.Range("H105").Value = 333334940
$ws.Range("I105").Value = 2400.5
$ws.Range("J105").Value = 1000000000
$ws.Range("K105").Value = 2400.5
$ws.Range("L105").Value = 1000000000
$ws.Range("M105").Value = -653.5
$ws.Range("N105").Value = -1000003494

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 9556138
$ws.Range("I113").Value = 6410855
$ws.Range("J113").Value = 12963528
$ws.Range("K113").Value = 19232565
$ws.Range("L113").Value = 38890584
$ws.Range("M113").Value = -19230395
$ws.Range("N113").Value = -38894924
$ws.Range("H131").Value = 811.08
$ws.Range("I131").Value = 571.5
$ws.Range("J131").Value = 826.3723
$ws.Range("K131").Value = 1714.5
$ws.Range("L131").Value = 2479.1169
$ws.Range("M131").Value = 3325.5
$ws.Range("N131").Value = -12559.1169

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 813.4483
$ws.Range("I97").Value = 915.94446
$ws.Range("J97").Value = 645.7273
$ws.Range("K97").Value = 915.94446
$ws.Range("L97").Value = 645.7273
$ws.Range("M97").Value = -419.94446
$ws.Range("N97").Value = -1637.7273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 990.1667
$ws.Range("I16").Value = 987.8
$ws.Range("J16").Value = 1002
$ws.Range("K16").Value = 987.8
$ws.Range("L16").Value = 1002
$ws.Range("M16").Value = -817.8
$ws.Range("N16").Value = -1342
$ws.Range("H40").Value = 22731010
$ws.Range("I40").Value = 3729.1428
$ws.Range("J40").Value = 62503750
$ws.Range("K40").Value = 3729.1428
$ws.Range("L40").Value = 62503750
$ws.Range("M40").Value = -3593.1428
$ws.Range("N40").Value = -62504022
$ws.Range("H61").Value = 26317712
$ws.Range("I61").Value = 1818
$ws.Range("J61").Value = 166669150
$ws.Range("K61").Value = 1818
$ws.Range("L61").Value = 166669150
$ws.Range("M61").Value = -1616
$ws.Range("N61").Value = -166669554
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H113").Value = 26317712
$ws.Range("I113").Value = 1818
$ws.Range("J113").Value = 166669150
$ws.Range("K113").Value = 1818
$ws.Range("L113").Value = 166669150
$ws.Range("M113").Value = 352
$ws.Range("N113").Value = -166673490

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 6570
$ws.Range("J41").Value = 6570
$ws.Range("L41").Value = 6570
$ws.Range("N41").Value = -7350
$ws.Range("H96").Value = 1840
$ws.Range("I96").Value = 1575
$ws.Range("J96").Value = 2900
$ws.Range("K96").Value = 1575
$ws.Range("L96").Value = 2900
$ws.Range("M96").Value = -202
$ws.Range("N96").Value = -5646
$ws.Range("H113").Value = 62503228
$ws.Range("I113").Value = 111116410
$ws.Range("J113").Value = 571.4286
$ws.Range("K113").Value = 333349230
$ws.Range("L113").Value = 1714.2858
$ws.Range("M113").Value = -333347060
$ws.Range("N113").Value = -6054.2858
